$wb = $excel.ActiveWorkbook

# --- sheet1 (SearchBarData): Adidas -> ADIDAS ---
$ws1 = $wb.Worksheets.Item("SearchBarData")
$ws1.Range("A1").Value = "ADIDAS"

# --- sheet2 (AssertData): adidas -> ADIDAS, add rows 8-15 ---
$ws2 = $wb.Worksheets.Item("AssertData")
$ws2.Range("A2").Value = "ADIDAS"
$ws2.Range("A8").Value = "BADMINTON RACKETS"
$ws2.Range("A9").Value = "YONEX ARCSABER 2 FEEL"
$ws2.Range("A10").Value = "CART"
$ws2.Range("A11").Value = "Product added to cart"
$ws2.Range("A12").Value = "YONEX"
$ws2.Range("A13").Value = "ZIPCODE"
$ws2.Range("A14").Value = "shoes"
$ws2.Range("A15").Value = "required"

# --- add new sheet InputData at the end of the workbook ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "InputData"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# re-fetch by name: the Move() call can leave the old handle pointing
# at the wrong (positionally reseated) sheet
$ws4 = $wb.Worksheets.Item("InputData")
$ws4.Range("A1").Value = "'641001"
